$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Instructions" sheet - update validation-rule text to reflect the new,
#    human-friendly column header names. The sheet is protected, so we need
#    to unprotect it first and re-protect it once we are done editing.
# ---------------------------------------------------------------------------
$wsInstructions = $wb.Worksheets.Item("Instructions")
$wsInstructions.Unprotect()

$wsInstructions.Range("A7").Value  = 'Please note that no field in the "Forecast Report" sheet may be blank'
$wsInstructions.Range("A10").Value = '(1) "Model Year" should be a 4 digit integer'
$wsInstructions.Range("A11").Value = '(2) "Make" should be no more than 250 characters'
$wsInstructions.Range("A12").Value = '(3) "Model" should be no more than 250 characters'
$wsInstructions.Range("A13").Value = '(4) "Type" should be exactly one of: BEV, PHEV, FCEV, EREV'
$wsInstructions.Range("A14").Value = '(5) "Range" should be a real number with no more than 2 decimal places'
$wsInstructions.Range("A15").Value = '(6) "ZEV Class" should be a single, uppercase letter'
$wsInstructions.Range("A16").Value = '(7) "Vehicle Class and Interior Volume" should be no more than 250 characters'
$wsInstructions.Range("A17").Value = '(8) "Total ZEVs Supplied" should be an integer'

$wsInstructions.Protect([System.Reflection.Missing]::Value, $true, $true, $true)

# ---------------------------------------------------------------------------
# 2) "Forecast Report" sheet - rename headers, bold them, widen a couple of
#    columns and add dropdown (list) data validation driven off a new,
#    hidden "Dropdowns" sheet.
# ---------------------------------------------------------------------------
$wsForecast = $wb.Worksheets.Item("Forecast Report")

$wsForecast.Range("A1").Value = "Model Year"
$wsForecast.Range("B1").Value = "Make"
$wsForecast.Range("C1").Value = "Model"
$wsForecast.Range("D1").Value = "Type"
$wsForecast.Range("E1").Value = "Range"
$wsForecast.Range("F1").Value = "ZEV Class"
$wsForecast.Range("G1").Value = "Vehicle Class and Interior Volume"
$wsForecast.Range("H1").Value = "Total ZEVs Supplied"

$wsForecast.Range("A1:H1").Font.Bold = $true

$wsForecast.Columns.Item(7).ColumnWidth = 28.5
$wsForecast.Columns.Item(8).ColumnWidth = 16.5

# ---------------------------------------------------------------------------
# 3) Create the new, hidden "Dropdowns" sheet (placed after "Forecast
#    Report") which backs the dropdown lists on the Forecast Report sheet.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsDropdowns = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsDropdowns.Name = "Dropdowns"

$wsDropdowns.Range("A1").Value = "Type"
$wsDropdowns.Range("B1").Value = "ZEV Class"
$wsDropdowns.Range("C1").Value = "Vehicle Class and Interior Volume"
$wsDropdowns.Range("A1:C1").Font.Bold = $true

$wsDropdowns.Range("A2").Value = "BEV"
$wsDropdowns.Range("A3").Value = "PHEV"
$wsDropdowns.Range("A4").Value = "EREV"
$wsDropdowns.Range("A5").Value = "FCEV"

$wsDropdowns.Range("B2").Value = "A"
$wsDropdowns.Range("B3").Value = "B"
$wsDropdowns.Range("B4").Value = "C"

$wsDropdowns.Range("C2").Value  = "Two-seater"
$wsDropdowns.Range("C3").Value  = "Minicompact (less than 85 cu. ft.)"
$wsDropdowns.Range("C4").Value  = "Subcompact (85–99 cu. ft.)"
$wsDropdowns.Range("C5").Value  = "Compact (100–109 cu. ft.)"
$wsDropdowns.Range("C6").Value  = "Mid-size (110–119 cu. ft.)"
$wsDropdowns.Range("C7").Value  = "Full-size (120 cu. ft. or more)"
$wsDropdowns.Range("C8").Value  = "Station wagon: Small (less than 130 cu. ft.)"
$wsDropdowns.Range("C9").Value  = "Station wagon: Mid-size (130–159 cu. ft.)"
$wsDropdowns.Range("C10").Value = "Pickup truck: Small (less than 2,722 kg)"
$wsDropdowns.Range("C11").Value = "Pickup truck: Standard (2,722–3,856 kg)"
$wsDropdowns.Range("C12").Value = "Sport utility vehicle: Small (less than 2,722 kg)"
$wsDropdowns.Range("C13").Value = "Sport utility vehicle: Standard (2,722–4,536 kg)"
$wsDropdowns.Range("C14").Value = "Minivan (less than 3,856 kg)"
$wsDropdowns.Range("C15").Value = "Van: Cargo (less than 3,856 kg)"
$wsDropdowns.Range("C16").Value = "Van: Passenger (less than 4,536 kg)"
$wsDropdowns.Range("C17").Value = "Special purpose vehicle (less than 3,856 kg)"
$wsDropdowns.Range("C18").Value = "Other/TBD"

# ---------------------------------------------------------------------------
# 4) Wire up the dropdown (list) data validations on "Forecast Report",
#    sourced from the "Dropdowns" sheet.
# ---------------------------------------------------------------------------
$wsForecast.Range("D2:D200").Validation.Add(
    [Microsoft.Office.Interop.Excel.XlDVType]::xlValidateList,
    [Microsoft.Office.Interop.Excel.XlDVAlertStyle]::xlValidAlertStop,
    [Microsoft.Office.Interop.Excel.XlFormatConditionOperator]::xlBetween,
    "=Dropdowns!`$A`$2:`$A`$5")

$wsForecast.Range("F2:F200").Validation.Add(
    [Microsoft.Office.Interop.Excel.XlDVType]::xlValidateList,
    [Microsoft.Office.Interop.Excel.XlDVAlertStyle]::xlValidAlertStop,
    [Microsoft.Office.Interop.Excel.XlFormatConditionOperator]::xlBetween,
    "=Dropdowns!`$B`$2:`$B`$4")

$wsForecast.Range("G2:G200").Validation.Add(
    [Microsoft.Office.Interop.Excel.XlDVType]::xlValidateList,
    [Microsoft.Office.Interop.Excel.XlDVAlertStyle]::xlValidAlertStop,
    [Microsoft.Office.Interop.Excel.XlFormatConditionOperator]::xlBetween,
    "=Dropdowns!`$C`$2:`$C`$18")

# ---------------------------------------------------------------------------
# 5) Zoom "Forecast Report" to 120%, hide "Dropdowns" and restore
#    "Instructions" as the selected/active sheet.
# ---------------------------------------------------------------------------
$wsForecast.Activate()
$excel.ActiveWindow.Zoom = 120

$wsDropdowns.Visible = $false

$wsInstructions.Activate()
